# Updates odds data on row 3 and row 4 of Sheet1 to reflect the latest
# FlashScore odds ("Atualizando o arquivo XLSX").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 3 updates ---
$ws.Range("G3").Value = 2.63
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 2.55
$ws.Range("J3").Value = 3.4
$ws.Range("K3").Value = 2.05
$ws.Range("L3").Value = 3.25

$ws.Range("W3").Value = 8.5
$ws.Range("X3").Value = 13
$ws.Range("Z3").Value = 29
$ws.Range("AA3").Value = 23
$ws.Range("AB3").Value = 34

$ws.Range("AH3").Value = 8
$ws.Range("AI3").Value = 12
$ws.Range("AJ3").Value = 10
$ws.Range("AK3").Value = 26
$ws.Range("AL3").Value = 21

$ws.Range("AN3").Value = 4.75
$ws.Range("AP3").Value = 26
$ws.Range("AR3").Value = 81

$ws.Range("AX3").Value = 4.5
$ws.Range("AY3").Value = 15
$ws.Range("BB3").Value = 67

# --- Row 4 updates ---
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 3.4
$ws.Range("Q4").Value = 2.03
$ws.Range("R4").Value = 1.83
